# Insert a new data row at row 562 (pushing existing rows 562..576 down to
# 563..577) and populate it with the new weekly price entry for Mandarina
# "Murcott" reported at the "Macroferia Regional de Talca" fair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 562:576 down by one row.
$ws.Rows.Item(562).Insert()

# Fill in the new row 562 with the new observation.
$ws.Cells.Item(562, 1).Value2  = 5
$ws.Cells.Item(562, 2).Value2  = 'Macroferia Regional de Talca'
$ws.Cells.Item(562, 3).Value2  = 'Maule'
$ws.Cells.Item(562, 4).Value2  = 45239
$ws.Cells.Item(562, 5).Value2  = 7
$ws.Cells.Item(562, 6).Value2  = 'Fruta'
$ws.Cells.Item(562, 7).Value2  = 100102
$ws.Cells.Item(562, 8).Value2  = 'Cítricos'
$ws.Cells.Item(562, 9).Value2  = 100102004
$ws.Cells.Item(562, 10).Value2 = 'Mandarina'
$ws.Cells.Item(562, 11).Value2 = 'Murcott'
$ws.Cells.Item(562, 12).Value2 = 'Primera'
$ws.Cells.Item(562, 13).Value2 = 400
$ws.Cells.Item(562, 14).Value2 = 8000
$ws.Cells.Item(562, 15).Value2 = 8000
$ws.Cells.Item(562, 16).Value2 = 8000
$ws.Cells.Item(562, 17).Value2 = '$/bandeja 18 kilos'
$ws.Cells.Item(562, 18).Value2 = 'Región de O''Higgins'
$ws.Cells.Item(562, 19).Value2 = 444
$ws.Cells.Item(562, 20).Value2 = 18
